# Inventory.xlsx - "Add supplier / Group logic"
# Populates supplier (E) / group (M) data for the first few inventory rows,
# tweaks a couple of numeric columns (Ціна за одиницю / Нижня межа поповнення),
# adds values to the previously-empty "Інтервал поповнення в днях" (J) column,
# and restores the UI selection to G6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InventoryList")

# --- Row 4 ---------------------------------------------------------------
$ws.Range("C4").Value = "Id01"
$ws.Range("D4").Value = "Tovar01"
$ws.Range("E4").Value = "Sup1"
$ws.Range("F4").Value = 40
$ws.Range("I4").Value = 32
$ws.Range("J4").Value = 20
$ws.Range("M4").Value = "Group2"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("C5").Value = "Id02"
$ws.Range("D5").Value = "Tovar02"
$ws.Range("E5").Value = "Sup2"
$ws.Range("F5").Value = 22
$ws.Range("I5").Value = 32
$ws.Range("J5").Value = 14
$ws.Range("M5").Value = "Group1"

# --- Row 6 ---------------------------------------------------------------
$ws.Range("C6").Value = "Id03"
$ws.Range("D6").Value = "Tovar03"
$ws.Range("E6").Value = "Sup2"
$ws.Range("F6").Value = 43
$ws.Range("I6").Value = 33
$ws.Range("J6").Value = 55
$ws.Range("M6").Value = "Group2"

# --- Restore the active selection on the sheet ---------------------------
$ws.Activate() | Out-Null
$ws.Range("G6").Select() | Out-Null
